$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 15
$ws.Range("C3").Value = 15
$ws.Range("C4").Value = 15

$ws.Range("C8").Value = 20
$ws.Range("C9").Value = 20
$ws.Range("C10").Value = 20
